$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, '516010', '2025-06-16', '游戏ETF', 1.25, 5.24, 24.87, 2.13, 1.181, 1.1574, 1.123, 0.007600000000000051, $true, 0.01549999999999985, $false),
    @(3, '512980', '2025-06-16', '传媒ETF', 0.84, 2.55, 9.04, 1.03, 0.828, 0.8196999999999999, 0.808, 0.001449999999999951, $false, 0.004399999999999848, $false),
    @(4, '515230', '2025-06-16', '软件ETF', 0.8, 2.43, 2.3, 0.71, 0.802, 0.7976, 0.792, -0.001299999999999968, $false, 0.003499999999999948, $false),
    @(5, '515880', '2025-06-16', '通信ETF', 1.35, 2.34, -0.81, 2.74, 1.32, 1.2867, 1.257, 0.004200000000000204, $true, 0.0132000000000001, $false),
    @(6, '512200', '2025-06-16', '房地产ETF', 1.36, 2.03, -7.05, 2.07, 1.357, 1.3526, 1.351, -0.000600000000000378, $false, 0.002799999999999914, $false),
    @(7, '516510', '2025-06-16', '云计算ETF', 1.13, 1.89, 2.72, 1.64, 1.141, 1.1268, 1.117, -0.001749999999999918, $false, 0.004499999999999948, $false),
    @(8, '512000', '2025-06-16', '券商ETF', 1.04, 1.46, -6.79, 8.99, 1.032, 1.021, 1.016, -0.0005500000000002725, $false, 0.004900000000000126, $false),
    @(9, '513330', '2025-06-16', '恒生互联网ETF', 0.49, 1.45, 17.75, 19, 0.494, 0.4837, 0.479, 0, $false, 0.002299999999999969, $false),
    @(10, '513290', '2025-06-16', '纳指生物科技ETF', 1.13, 1.35, -2.42, 0.5, 1.122, 1.1109, 1.097, 0.003049999999999997, $false, 0.003900000000000015, $false),
    @(11, '512800', '2025-06-16', '银行ETF', 1.67, 1.21, 12.46, 3.66, 1.652, 1.6399, 1.624, 0.002550000000000052, $true, 0.006500000000000172, $false),
    @(12, '510230', '2025-06-16', '金融ETF', 1.4, 1.16, 7.12, 0.28, 1.383, 1.3721, 1.362, 0.001299999999999857, $true, 0.005400000000000071, $false),
    @(13, '513520', '2025-06-16', '日经ETF', 1.48, 1.16, 1.58, 0.76, 1.468, 1.4647, 1.457, 0.001500000000000057, $true, 0.002500000000000169, $false),
    @(14, '510900', '2025-06-16', 'H股ETF', 1.12, 0.99, 18.42, 2.55, 1.127, 1.1136, 1.107, 0.0003999999999999559, $true, 0.003100000000000103, $false),
    @(15, '159949', '2025-06-16', '创业板50ETF', 0.92, 0.99, -4.79, 5.18, 0.909, 0.8999, 0.9, -0.000300000000000078, $false, 0.004099999999999993, $false),
    @(16, '515070', '2025-06-16', '人工智能AIETF', 1.17, 0.95, 0.86, 0.76, 1.181, 1.1689, 1.164, -0.001349999999999962, $false, 0.003500000000000059, $false),
    @(17, '159915', '2025-06-16', '创业板ETF', 2.04, 0.84, -3.09, 10.11, 2.028, 2.0102, 2.004, -0.000700000000000145, $false, 0.008100000000000218, $false),
    @(18, '513100', '2025-06-16', '纳指ETF', 1.57, 0.84, -3.81, 7.23, 1.571, 1.5709, 1.561, 0.000199999999999978, $true, 0.000199999999999978, $false),
    @(19, '515260', '2025-06-16', '电子ETF', 0.85, 0.83, -6.17, 0.06, 0.854, 0.8474999999999999, 0.846, -0.001049999999999773, $false, 0.002299999999999969, $false),
    @(20, '513500', '2025-06-16', '标普500ETF', 2.04, 0.74, -5.74, 2.33, 2.036, 2.0281, 2.015, 0.001349999999999962, $true, 0.002799999999999692, $false),
    @(21, '513800', '2025-06-16', '日本东证指数ETF', 1.48, 0.68, 8.12, 0.11, 1.477, 1.4811, 1.472, 0.00154999999999994, $true, 0.000500000000000167, $false),
    @(22, '159781', '2025-06-16', '科创创业ETF', 0.53, 0.56, -5.49, 0.38, 0.535, 0.5321, 0.531, -0.0007000000000000339, $false, 0.0009000000000000119, $false),
    @(23, '512890', '2025-06-16', '红利低波ETF', 1.18, 0.51, 4.54, 2.73, 1.172, 1.1678, 1.158, 0.001400000000000068, $true, 0.002600000000000158, $true),
    @(24, '515790', '2025-06-16', '光伏ETF', 0.65, 0.47, -14.66, 1.84, 0.651, 0.6476, 0.653, -0.001550000000000051, $false, 0.0007999999999999119, $false),
    @(25, '510760', '2025-06-16', '上证综指ETF', 1.09, 0.46, 3.32, 0.36, 1.085, 1.0799, 1.075, 0.0004999999999999449, $true, 0.002699999999999925, $false),
    @(26, '159667', '2025-06-16', '工业母机ETF', 1.11, 0.45, 9.99, 0.32, 1.118, 1.1119, 1.116, -0.001850000000000129, $false, 0.001499999999999835, $false),
    @(27, '512500', '2025-06-16', '中证500ETF华夏', 3.18, 0.44, 1.4, 4.01, 3.182, 3.1632, 3.145, -0.00004999999999988347, $false, 0.007899999999999796, $false),
    @(28, '510210', '2025-06-16', '上证指数ETF', 0.84, 0.36, 3.43, 0.98, 0.842, 0.8390000000000001, 0.834, 0.000400000000000178, $true, 0.001600000000000046, $false),
    @(29, '515250', '2025-06-16', '智能汽车ETF', 0.93, 0.32, -0.43, 0.25, 0.946, 0.9488999999999999, 0.952, -0.002049999999999996, $false, -0.000700000000000145, $false),
    @(30, '512480', '2025-06-16', '半导体ETF', 1, 0.3, -1.29, 5.16, 1.014, 1.0115, 1.013, -0.002350000000000074, $false, 0.001000000000000112, $false),
    @(31, '515800', '2025-06-16', '800ETF', 1.01, 0.3, -0.2, 0.34, 1.01, 1.006, 1.006, -0.0004999999999999449, $false, 0.001599999999999824, $false),
    @(32, '159691', '2025-06-16', '港股红利ETF', 1.21, 0.25, 8.54, 3.06, 1.194, 1.1809, 1.165, 0.003549999999999942, $true, 0.00529999999999986, $true),
    @(33, '510050', '2025-06-16', '上证50ETF', 2.75, 0.25, 0.55, 15.9, 2.753, 2.7526, 2.764, -0.003349999999999742, $false, 0.000300000000000189, $false),
    @(34, '510300', '2025-06-16', '沪深300ETF', 3.99, 0.25, -0.8, 25.87, 3.989, 3.9787, 3.983, -0.002750000000000252, $false, 0.00480000000000036, $false),
    @(35, '159770', '2025-06-16', '机器人ETF', 0.85, 0.24, 5.58, 0.95, 0.864, 0.8648999999999999, 0.876, -0.003449999999999953, $false, -0.001400000000000179, $false),
    @(36, '512760', '2025-06-16', '芯片ETF', 1.1, 0.18, -2.13, 1.28, 1.12, 1.1189, 1.123, -0.003149999999999986, $false, 0.00009999999999998899, $false),
    @(37, '512690', '2025-06-16', '酒ETF', 0.55, 0.18, -11.36, 9.58, 0.568, 0.574, 0.583, -0.003149999999999875, $false, -0.003700000000000037, $false),
    @(38, '515220', '2025-06-16', '煤炭ETF', 0.99, 0.1, -12.97, 2.04, 0.99, 0.9878, 0.989, -0.00005000000000010552, $false, 0.000400000000000178, $false),
    @(39, '510410', '2025-06-16', '资源ETF', 1.26, 0.08, 3.54, 0.09, 1.236, 1.2237, 1.217, 0.002250000000000085, $true, 0.005599999999999827, $false),
    @(40, '515900', '2025-06-16', '央企创新驱动ETF', 1.43, 0.07, -4.98, 0.21, 1.432, 1.4269, 1.426, -0.0007999999999999119, $false, 0.001700000000000035, $false),
    @(41, '515080', '2025-06-16', '中证红利ETF', 1.53, 0.07, -1.14, 1.5, 1.53, 1.5283, 1.527, -0.000199999999999978, $false, 0.001300000000000079, $false),
    @(42, '511090', '2025-06-16', '30年国债ETF', 124.26, 0.01, 1.31, 52.32, 123.889, 123.4057, 123.301, 0.05725000000002467, $false, 0.1278999999999968, $false),
    @(43, '512660', '2025-06-16', '军工ETF', 1.06, 0, 1.34, 4.58, 1.053, 1.0505, 1.044, 0.0001500000000000945, $false, 0.00340000000000007, $false),
    @(44, '159666', '2025-06-16', '交通运输ETF', 0.99, 0, -0.5, 0.02, 0.986, 0.9865999999999999, 0.985, -0.000400000000000178, $false, 0.000500000000000056, $false),
    @(45, '561560', '2025-06-16', '电力ETF', 1.17, -0.09, -2.09, 0.32, 1.172, 1.1724, 1.178, -0.0003999999999999559, $false, -0.0007999999999999119, $false),
    @(46, '560070', '2025-06-16', '央企红利ETF基金', 1.02, -0.1, -3.87, 0.02, 1.023, 1.0236, 1.024, -0.0005500000000000504, $false, -0.0000000000000002220446049250313, $false),
    @(47, '588000', '2025-06-16', '科创50ETF', 1.02, -0.1, -2.3, 14.45, 1.037, 1.0378, 1.038, -0.002249999999999863, $false, -0.000299999999999967, $false),
    @(48, '515210', '2025-06-16', '钢铁ETF', 1.2, -0.17, 2.13, 0.45, 1.205, 1.2013, 1.21, -0.001799999999999802, $false, 0.000199999999999978, $false),
    @(49, '516020', '2025-06-16', '化工ETF', 0.6, -0.17, -1.15, 0.09, 0.598, 0.5932, 0.593, -0.000200000000000089, $false, 0.001599999999999935, $false),
    @(50, '518880', '2025-06-16', '黄金ETF', 7.56, -0.33, 27.58, 50.41, 7.455, 7.4378, 7.378, 0.01414999999999988, $false, 0.01850000000000041, $false),
    @(51, '159928', '2025-06-16', '消费ETF', 0.8, -0.38, -2.69, 2.21, 0.816, 0.8187999999999999, 0.82, -0.001650000000000151, $false, -0.002000000000000113, $false),
    @(52, '516670', '2025-06-16', '畜牧养殖ETF', 0.66, -0.45, 6.6, 0.2, 0.667, 0.6637000000000001, 0.657, 0.0005499999999999394, $true, 0.001600000000000046, $false),
    @(53, '159637', '2025-06-16', '新能源车龙头ETF', 0.59, -0.5, 0.68, 0.11, 0.598, 0.5984, 0.607, -0.001599999999999935, $false, -0.0007999999999999119, $false),
    @(54, '159652', '2025-06-16', '有色50ETF', 0.96, -0.52, 12.66, 0.15, 0.944, 0.9315, 0.926, 0.00154999999999994, $true, 0.004999999999999893, $false),
    @(55, '159825', '2025-06-16', '农业ETF', 0.72, -0.55, 7.96, 0.45, 0.721, 0.7127, 0.7, 0.00154999999999994, $true, 0.003599999999999937, $false),
    @(56, '512170', '2025-06-16', '医疗ETF', 0.33, -0.6, 0.61, 4.61, 0.333, 0.3319, 0.328, 0.0001499999999999835, $true, 0.0007999999999999674, $false),
    @(57, '562860', '2025-06-16', '生物疫苗ETF', 0.66, -0.61, 7.73, 0.1, 0.662, 0.6562, 0.636, 0.002900000000000014, $true, 0.004199999999999982, $false),
    @(58, '562390', '2025-06-16', '中药50ETF', 0.97, -0.61, -2.8, 0.01, 0.985, 0.9833000000000001, 0.974, 0.0006999999999999229, $true, 0.000700000000000145, $true),
    @(59, '159643', '2025-06-16', '疫苗ETF', 0.58, -0.68, -1.35, 0.06, 0.59, 0.5871, 0.577, 0.0009500000000000064, $true, 0.001800000000000024, $true),
    @(60, '512010', '2025-06-16', '医药ETF', 0.37, -0.8, 3.6, 7.81, 0.376, 0.3753, 0.371, 0.000350000000000017, $true, 0.0003000000000000225, $true),
    @(61, '159883', '2025-06-16', '医疗器械ETF', 0.48, -0.83, -1.65, 0.26, 0.488, 0.4866, 0.48, 0.0001000000000000445, $true, 0.0007000000000000894, $false),
    @(62, '513120', '2025-06-16', '港股创新药ETF', 1.18, -0.84, 62.76, 81.08, 1.146, 1.0922, 1.02, 0.01469999999999994, $true, 0.02049999999999996, $false),
    @(63, '513060', '2025-06-16', '恒生医疗ETF', 0.58, -0.86, 48.46, 26.28, 0.567, 0.5463, 0.52, 0.005200000000000093, $true, 0.007900000000000018, $false)
)

# Force column B (date) to be treated as text so date-like strings
# ("2025-06-16") are not auto-converted to date serials.
$ws.Range("B2:B63").NumberFormat = "@"

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
}

# Restore the General format now that the text values are committed -
# Excel only re-parses text into a date on *entry*, not when the display
# format changes, so the strings stay intact.
$ws.Range("B2:B63").NumberFormat = "General"
$ws.Range("B2:B63").ClearFormats()
